$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 154 (pushing existing data down), then
# resize the "Snippets" table to include them.
$ws.Rows("154:155").Insert()

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F188"))

# Populate the two new rows with the "Setting" class rows (key/value props).
$ws.Range("A154").Value = "Word"
$ws.Range("B154").Value = "Setting"
$ws.Range("C154").Value = "key"
$ws.Range("E154").Value = "word-document-manage-settings"
$ws.Range("F154").Value = "addEditSetting"

$ws.Range("A155").Value = "Word"
$ws.Range("B155").Value = "Setting"
$ws.Range("C155").Value = "value"
$ws.Range("E155").Value = "word-document-manage-settings"
$ws.Range("F155").Value = "addEditSetting"

# Restore the view state (scrolled/selected differently after the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("C156").Select()
